# "Add files via upload" — the uploaded workbook drops the bottom
# "Seri Açıklamaları / Notlar" metadata block (EVDS export footer) that used
# to live in rows 173-190, while leaving the two bold section-header cells
# (A173, A177) in place but empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fully clear (contents + formatting) the metadata rows that disappear
# entirely from the sheet.
$ws.Range("A174:C176").Clear()
$ws.Range("A178:C190").Clear()

# The two section-header cells stay on the sheet (they keep their bold
# style, s="1") but lose their text.
$ws.Range("A173").ClearContents()
$ws.Range("A177").ClearContents()

# Match the saved selection/view state left after trimming the sheet.
$ws.Range("H195").Select() | Out-Null
